# Apply cryptos list update (prices & 1h volume changes) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "29.199.09"
$ws.Cells.Item(2, 5).Value = "  -0.63%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.859.28"
$ws.Cells.Item(3, 5).Value = "  -1.26%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.02%  "

# Row 5
$ws.Cells.Item(5, 2).Value = "XRP"
$ws.Cells.Item(5, 3).Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "0.7035"
$ws.Cells.Item(5, 5).Value = "  -1.30%  "

# Row 6
$ws.Cells.Item(6, 2).Value = "BNB"
$ws.Cells.Item(6, 3).Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "242.22"
$ws.Cells.Item(6, 5).Value = "  -0.07%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.02%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3111"
$ws.Cells.Item(8, 5).Value = "  -0.47%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.07780"
$ws.Cells.Item(9, 5).Value = "  -3.26%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "24.19"
$ws.Cells.Item(10, 5).Value = "  -4.45%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -3.98%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "1.861.70"
$ws.Cells.Item(12, 5).Value = "  -2.39%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "5.174"
$ws.Cells.Item(13, 5).Value = "  -1.40%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "93.41"
$ws.Cells.Item(14, 5).Value = "  -0.19%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.6957"
$ws.Cells.Item(15, 5).Value = "  -3.37%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "6.351"
$ws.Cells.Item(16, 5).Value = "  +0.25%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "29.185.41"
$ws.Cells.Item(17, 5).Value = "  -0.71%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.000008286"
$ws.Cells.Item(18, 5).Value = "  -2.70%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "251.00"
$ws.Cells.Item(19, 5).Value = "  +3.89%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "2.116.12"
$ws.Cells.Item(20, 5).Value = "  -1.38%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  -1.31%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "1.001"
$ws.Cells.Item(22, 5).Value = "  -0.05%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "7.506"
$ws.Cells.Item(23, 5).Value = "  -4.45%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "1.000"
$ws.Cells.Item(24, 5).Value = "  -0.07%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.1549"
$ws.Cells.Item(25, 5).Value = "  -2.28%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "8.964"

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "159.34"
$ws.Cells.Item(27, 5).Value = "  -2.87%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "18.77"
$ws.Cells.Item(28, 5).Value = "  +1.02%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  -0.88%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "4.286"
$ws.Cells.Item(30, 5).Value = "  -3.02%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "4.256"
$ws.Cells.Item(31, 5).Value = "  -2.09%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  +0.75%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.05248"
$ws.Cells.Item(33, 5).Value = "  -2.23%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  -3.67%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.7423"
$ws.Cells.Item(35, 5).Value = "  -1.05%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -2.47%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.711"
$ws.Cells.Item(37, 5).Value = "  +0.45%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.01864"
$ws.Cells.Item(38, 5).Value = "  -1.37%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "1.244.33"
$ws.Cells.Item(39, 5).Value = "  -3.47%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.736"
$ws.Cells.Item(40, 5).Value = "  -0.21%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "6.222"
$ws.Cells.Item(41, 5).Value = "  -5.97%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "110.77"
$ws.Cells.Item(42, 5).Value = "  -1.03%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.8946"
$ws.Cells.Item(43, 5).Value = "  -3.13%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "71.20"
$ws.Cells.Item(44, 5).Value = "  -4.22%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.000"

# Row 46
$ws.Cells.Item(46, 5).Value = "  -0.36%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "2.012.39"
$ws.Cells.Item(47, 5).Value = "  -1.69%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.5182"
$ws.Cells.Item(48, 5).Value = "  -0.72%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "1.779"
$ws.Cells.Item(49, 5).Value = "  -1.55%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "9.407"
$ws.Cells.Item(50, 5).Value = "  -1.16%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.4292"
$ws.Cells.Item(51, 5).Value = "  -2.29%  "
